$wb = $excel.ActiveWorkbook

# Sheet references (tab order: Sheet1, 解說, Sheet3)
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 解說 (Sheet2) explanation text: TRUE/FALSE -> M/F wording ---
$ws2.Range("C3").Value = "男 - M "
$ws2.Range("D3").Value = "女 - F "

# --- Sheet1 gender cells: boolean TRUE/FALSE -> string "M"/"F" ---
$ws1.Range("E2").Value = "M"
$ws1.Range("E3").Value = "F"

# --- Update selections shown in each sheet, keeping Sheet1 as the active tab ---
$ws2.Range("F3").Select() | Out-Null
$ws1.Range("E5").Select() | Out-Null
